$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Jengibre at Terminal La Palmera
# de La Serena. It becomes the new row 152; every existing row from 152
# down to 172 shifts down by one (to 153..173) keeping its original data.
$ws.Rows.Item(152).Insert()

$ws.Range("A152").Value = 8
$ws.Range("B152").Value = "Terminal La Palmera de La Serena"
$ws.Range("C152").Value = "Coquimbo"
$ws.Range("D152").Value = 45154
$ws.Range("E152").Value = 4
$ws.Range("F152").Value = 100114007
$ws.Range("G152").Value = "Jengibre"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 520
$ws.Range("K152").Value = 18000
$ws.Range("L152").Value = 19000
$ws.Range("M152").Value = 18500
$ws.Range("N152").Value = "$/caja 13 kilos"
$ws.Range("O152").Value = "Perú"
$ws.Range("P152").Value = 1423
$ws.Range("Q152").Value = 13
$ws.Range("R152").Value = "Hortaliza"
